# Add two new arrival rows (14 and 15) for Friday, Jan 13 flights,
# mirroring the style/format of the existing last row (row 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 14) into the
# two new rows so styles/column layout (incl. blank K/M cells) stay intact.
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(15).PasteSpecial()
$ws.Rows.Item(16).PasteSpecial()

# Row 15 (arrival #14): Lauda Europe flight from London (STN)
$ws.Range("A15").Value = 14.0
$ws.Range("B15").Value = "Friday, Jan 13"
$ws.Range("C15").Value = "5:40 PM"
$ws.Range("D15").Value = "FR2468"
$ws.Range("E15").Value = "London"
$ws.Range("F15").Value = "(STN)"
$ws.Range("G15").Value = "Lauda Europe "
$ws.Range("H15").Value = "A320"
$ws.Range("I15").Value = "(9H-LOA)"
$ws.Range("J15").Value = "6:48 PM"
$ws.Range("K15").Borders.LineStyle = 0
$ws.Range("L15").Value = "1 hours, 8 minutes"
$ws.Range("M15").Borders.LineStyle = 0

# Row 16 (arrival #15): Ryanair flight from Dublin (DUB)
$ws.Range("A16").Value = 15.0
$ws.Range("B16").Value = "Friday, Jan 13"
$ws.Range("C16").Value = "7:00 PM"
$ws.Range("D16").Value = "FR1978"
$ws.Range("E16").Value = "Dublin"
$ws.Range("F16").Value = "(DUB)"
$ws.Range("G16").Value = "Ryanair "
$ws.Range("H16").Value = "B738"
$ws.Range("I16").Value = "(EI-DHZ)"
$ws.Range("J16").Value = "6:37 PM"
$ws.Range("K16").Borders.LineStyle = 0
$ws.Range("L16").Value = "0 hours, -23 minutes"
$ws.Range("M16").Borders.LineStyle = 0
